$d = $word.ActiveDocument

function Find-ParaIndexExact($doc, $text) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $pt = $p.Range.Text
        $pt = $pt.TrimEnd([char]13)
        if ($pt -eq $text) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndexContains($doc, $substr) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Reorder the three summary lines:
#      Aggregation / Activation / Alignment  ->  Activation / Alignment / Aggregation
#    and fix the wording "Contexts" -> "Context" in each, by editing the text
#    of each paragraph in place (keeps paragraph/run formatting untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Aggregation: Clustering (Contexts types Occurrences).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Activation: Classification (Context types Occurrences Attributes).", 2) | Out-Null

$d.Content.Find.Execute(
    "Activation: Classification (Contexts types Occurrences Attributes).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Alignment: Regression (Context types Occurrences Attributes Values).", 2) | Out-Null

$d.Content.Find.Execute(
    "Alignment: Regression (Contexts types Occurrences Attributes Values).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aggregation: Clustering (Context types Occurrences).", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Insert the new "Augmentations" block right after the (now reordered)
#    "Aggregation: Clustering (Context types Occurrences)." line, but before
#    the existing blank paragraph that separates it from "Model Semantics:".
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexExact $d "Aggregation: Clustering (Context types Occurrences)."

$newLines = @(
    "",
    "Augmentations:",
    "",
    "(Context, Occurrence) : Value;",
    "",
    "Activation:",
    "(Statement, Resource) : Kind;",
    "",
    "Alignment:",
    "(Kind, Statement) : Resource;",
    "",
    "Aggregation:",
    "(Resource, Kind) : Statement;"
)

$curIdx = $idx
foreach ($line in $newLines) {
    $curPara = $d.Paragraphs.Item($curIdx)
    $r = $curPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $curIdx = $curIdx + 1
    if ($line -ne "") {
        $newPara = $d.Paragraphs.Item($curIdx)
        $newPara.Range.Text = $line
    }
}

# ---------------------------------------------------------------------------
# 3) Remove the stray trailing empty run in the "Dimensions Encoding: ..."
#    paragraph, leaving just the single run with the text.
# ---------------------------------------------------------------------------
$dimIdx = Find-ParaIndexContains $d "Dimensions Encoding: Given Dimensional Contexts"
$dimPara = $d.Paragraphs.Item($dimIdx)
$dimRange = $dimPara.Range
$dimStart = $dimRange.Start
$dimEnd = $dimRange.End
$dimTextRange = $d.Range($dimStart, $dimEnd - 1)
$dimTextRange.Delete()

$dimPara2 = $d.Paragraphs.Item($dimIdx)
$dimPara2.Range.Text = "Dimensions Encoding: Given Dimensional Contexts (CSPO Models set layouts) having four dimensional sets (Types Model, Individuals Model, Mappings Model, State Model) each representing (nested) CSPO inputs / parts of a recursively aggregated CSPO layout (i.e. aggregated layout Context is Mappings Model, Subject is State Model, etc.) having this setting (Models types / layers class / instance IDs) reified in this fifth " + [char]34 + "Focus" + [char]34 + " Model which represents a " + [char]34 + "snapshot" + [char]34 + " of current state and available transitions (Focus shifts)."

Write-Output "done"
